# Updated cryptos list on Tue Jun 27 03:34:10 UTC 2023 with GitHub Actions
#
# D-column price strings that look like plain numbers (single decimal point,
# e.g. "1.000", "237.89") get a leading apostrophe so Excel stores them as
# text instead of coercing them to a numeric value (which would silently
# drop trailing zeros / change precision). D-column strings that already
# contain multiple dots (e.g. "30.409.15") are never number-like, so a plain
# assignment is enough and keeps them as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.409.15"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.873.84"
$ws.Range("E3").Value = "  -0.21%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'237.89"
$ws.Range("E5").Value = "  +0.75%  "

# Row 6 - USDC
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  +0.03%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4820"
$ws.Range("E7").Value = "  -0.30%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.2826"
$ws.Range("E8").Value = "  -1.68%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.06555"
$ws.Range("E9").Value = "  -0.52%  "

# Row 10 - WrappedEther
$ws.Range("D10").Value = "1.882.75"
$ws.Range("E10").Value = "  +0.29%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.07460"
$ws.Range("E11").Value = "  +2.51%  "

# Row 12 - Solana
$ws.Range("D12").Value = "'16.40"
$ws.Range("E12").Value = "  -1.98%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'5.088"
$ws.Range("E13").Value = "  -1.56%  "

# Row 14 - Litecoin
$ws.Range("D14").Value = "'88.25"
$ws.Range("E14").Value = "  +1.36%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.6574"
$ws.Range("E15").Value = "  +0.45%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "30.333.42"
$ws.Range("E16").Value = "  +0.43%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "'13.30"
$ws.Range("E17").Value = "  -0.24%  "

# Row 18 - Dai
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  +0.06%  "

# Row 19 - ShibaInu (only Volume(1h) changes)
$ws.Range("E19").Value = "  -0.92%  "

# Row 20 - WrappedliquidstakedEther2.0
$ws.Range("D20").Value = "2.108.02"
$ws.Range("E20").Value = "  +0.49%  "

# Row 21 - now Uniswap (was BitcoinCash)
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.322"
$ws.Range("E21").Value = "  +0.22%  "

# Row 22 - now BinanceUSD (was Uniswap)
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.03%  "

# Row 23 - now BitcoinCash (was BinanceUSD)
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'221.34"
$ws.Range("E23").Value = "  +12.64%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "'6.184"
$ws.Range("E24").Value = "  +1.10%  "

# Row 25 - Cosmos
$ws.Range("D25").Value = "'9.262"
$ws.Range("E25").Value = "  -0.55%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'165.30"
$ws.Range("E26").Value = "  +3.81%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'18.61"
$ws.Range("E27").Value = "  +2.87%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "'1.961"
$ws.Range("E28").Value = "  +2.51%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "'1.455"
$ws.Range("E29").Value = "  +1.00%  "

# Row 30 - Stellar
$ws.Range("D30").Value = "'0.09379"
$ws.Range("E30").Value = "  +2.81%  "

# Row 31 - InternetComputer(DFINITY) (only Volume(1h) changes)
$ws.Range("E31").Value = "  +0.75%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'4.025"
$ws.Range("E32").Value = "  -0.72%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.05054"
$ws.Range("E33").Value = "  -1.37%  "

# Row 34 - ARBITRUM
$ws.Range("D34").Value = "'1.211"
$ws.Range("E34").Value = "  +10.58%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value = "'0.7562"
$ws.Range("E35").Value = "  +4.99%  "

# Row 36 - HuobiToken
$ws.Range("D36").Value = "'2.716"
$ws.Range("E36").Value = "  +0.20%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.01837"
$ws.Range("E37").Value = "  +2.22%  "

# Row 38 - MXToken (only Volume(1h) changes)
$ws.Range("E38").Value = "  -0.57%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "'2.084"
$ws.Range("E39").Value = "  +2.20%  "

# Row 40 - TrustWalletToken
$ws.Range("D40").Value = "'0.9052"
$ws.Range("E40").Value = "  -1.30%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "'5.940"
$ws.Range("E41").Value = "  +2.57%  "

# Row 42 - Quant
$ws.Range("D42").Value = "'106.77"
$ws.Range("E42").Value = "  +0.64%  "

# Row 43 - TheSandbox
$ws.Range("D43").Value = "'0.4294"
$ws.Range("E43").Value = "  +0.18%  "

# Row 44 - PaxDollar (only Volume(1h) changes)
$ws.Range("E44").Value = "  +0.49%  "

# Row 45 - Aptos
$ws.Range("D45").Value = "'7.467"
$ws.Range("E45").Value = "  +0.84%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'66.25"
$ws.Range("E46").Value = "  -0.29%  "

# Row 47 - Algorand (only Volume(1h) changes)
$ws.Range("E47").Value = "  -1.35%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "'1.483"
$ws.Range("E48").Value = "  +8.48%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "'8.933"
$ws.Range("E49").Value = "  -2.56%  "

# Row 50 - Elrond
$ws.Range("D50").Value = "'34.25"
$ws.Range("E50").Value = "  +0.83%  "

# Row 51 - Decentraland
$ws.Range("D51").Value = "'0.3895"
$ws.Range("E51").Value = "  +1.73%  "
